# Sync up for main test flow
$wb = $excel.ActiveWorkbook

# Remove the second worksheet ("1_8") - only "3_3" remains
$wb.Worksheets.Item("1_8").Delete()

$ws = $wb.Worksheets.Item("3_3")

# Update existing rows (case_0 / case_1) with new measured values
$ws.Range("B2").Value = -0.0044
$ws.Range("C2").Value = -0.0004
$ws.Range("D2").Value = -0.0104

$ws.Range("B3").Value = -0.0052
$ws.Range("C3").Value = 0.0009
$ws.Range("D3").Value = -0.0117
$ws.Range("E3").Value = 0.004

# New row 4: case_2 label (format copied from the row-3 label cell)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "case_2"

# New row 5: case_3 label (format copied from the row-3 label cell)
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "case_3"

# Raw-data strings for column G (order matches shared-string table layout)
$ws.Range("G2").Value = "-0.004065,-0.004718,-0.010383,-0.001952,-0.000356,-0.006857,-0.006769,-0.002027,-0.003243,-0.003141"
$ws.Range("G3").Value = "0.000942,0.000575,-0.007054,-0.009377,-0.00518,-0.011694,-0.003331,-0.007693,-0.00166,-0.007326"
$ws.Range("G4").Value = "-0.005471,-0.004432,-0.003522,0.000656,-0.006674,-0.0063,-0.004167,0.000853,-0.005661,-0.00092"
$ws.Range("G5").Value = "-0.00611,-0.006674,-0.007598,-0.009269,-0.002965,-0.007788,0.002137,-0.001477,-0.005465,-0.000825"

# Remaining numeric cells for the new rows
$ws.Range("B4").Value = -0.0036
$ws.Range("C4").Value = 0.0009
$ws.Range("D4").Value = -0.0067
$ws.Range("E4").Value = 0.0027
$ws.Range("F4").Value = 10

$ws.Range("B5").Value = -0.0046
$ws.Range("C5").Value = 0.0021
$ws.Range("D5").Value = -0.0093
$ws.Range("E5").Value = 0.0035
$ws.Range("F5").Value = 10
